# CobaltUsers.xlsx update
# - Rename Sheet2 -> Emails, add Email/Password rows
# - Sheet3: add "Y" flag value
# - Users sheet: add 28 new test users (rows 53-80) with hyperlinked emails

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet3: single flag cell used elsewhere in the workbook
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Range("A2").Value = "Y"

# ---------------------------------------------------------------------------
# Sheet2 -> Emails: small lookup table of email/password pairs
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Name = "Emails"

$ws2.Range("A1").Value = "Email"
$ws2.Range("B1").Value = "Password"
$ws2.Range("A2").Value = "tr-anz-tester1@yandex.com"
$ws2.Range("B2").Value = "tranztest"
$ws2.Range("A3").Value = "tr-anz-tester2@yandex.com"
$ws2.Range("B3").Value = "tranztest"

$ws2.Columns.Item(1).ColumnWidth = 26.28515625
$ws2.Columns.Item(2).ColumnWidth = 13.85546875
$ws2.Range("A1:B3").Select()

# ---------------------------------------------------------------------------
# Users sheet: append new test-user rows (53-80)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Users")

$newUsers = @(
    @("SearchOpenWebUser1", "SearchOpenWeb@mailinator.com"),
    @("FFHUser1", "FFHUser1@mailinator.com"),
    @("FFHUser2", "FFHUser2@mailinator.com"),
    @("FFHUser3", "FFHUser3@mailinator.com"),
    @("FFHUser4", "FFHUser4@mailinator.com"),
    @("FrontEndUser1", "FrontEndUser1@mailinator.com"),
    @("FrontEndUser2", "FrontEndUser2@mailinator.com"),
    @("FrontEndUser3", "FrontEndUser3@mailinator.com"),
    @("FrontEndUser4", "FrontEndUser4@mailinator.com"),
    @("FrontEndUser5", "FrontEndUser5@mailinator.com"),
    @("FrontEndUser6", "FrontEndUser6@mailinator.com"),
    @("FrontEndUser7", "FrontEndUser7@mailinator.com"),
    @("FrontEndUser8", "FrontEndUser8@mailinator.com"),
    @("FrontEndUser9", "FrontEndUser9@mailinator.com"),
    @("FrontEndUser10", "FrontEndUser10@mailinator.com"),
    @("UrlUser1", "UrlUser1@mailinator.com"),
    @("UrlUser2", "UrlUser2@mailinator.com"),
    @("UrlUser3", "UrlUser3@mailinator.com"),
    @("LinkingUser1", "LinkingUser1@mailinator.com"),
    @("LoginUser1", "LoginUser1@mailinator.com"),
    @("LoginUser2", "LoginUser2@mailinator.com"),
    @("LoginUser3", "LoginUser3@mailinator.com"),
    @("LoginUser4", "LoginUser4@mailinator.com"),
    @("LoginUser5", "LoginUser5@mailinator.com"),
    @("LoginUser6", "LoginUser6@mailinator.com"),
    @("LoginUser7", "LoginUser7@mailinator.com"),
    @("CpetUser1", "CpetUser1@mailinator.com"),
    @("CpetUser2", "CpetUser2@mailinator.com")
)

$row = 53
foreach ($pair in $newUsers) {
    $userName = $pair[0]
    $email = $pair[1]

    $ws1.Cells.Item($row, 1).Value = $userName
    $ws1.Cells.Item($row, 2).Value = "Password1"

    $eCell = $ws1.Cells.Item($row, 5)
    $eCell.Value = "THIS IS IN USE 24/7 - DO NOT USE!"
    $eCell.Borders.Item(7).LineStyle = 1
    $eCell.Borders.Item(10).LineStyle = 1

    $fCell = $ws1.Cells.Item($row, 6)
    $fCell.Value = "N"
    $fCell.Borders.Item(7).LineStyle = 1
    $fCell.Borders.Item(10).LineStyle = 1

    $gCell = $ws1.Cells.Item($row, 7)
    $gCell.Value = $email
    $ws1.Hyperlinks.Add($gCell, "mailto:$email")

    $row = $row + 1
}

# trailing bordered-but-empty rows (81-84), matching the block formatting
for ($r = 81; $r -le 84; $r++) {
    $eCell = $ws1.Cells.Item($r, 5)
    $eCell.Borders.Item(7).LineStyle = 1
    $eCell.Borders.Item(10).LineStyle = 1
}

$ws1.Columns.Item(1).AutoFit()
$ws1.Columns.Item(5).AutoFit()

$ws1.Range("C80").Select()
